$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": add column AA (10-jul) ---
$wsSpot = $wb.Worksheets.Item("Prix Spot")

$wsSpot.Range("Z1").Copy()
$wsSpot.Range("AA1").PasteSpecial(-4122)
$wsSpot.Range("AA1").Value = "10-jul"
$excel.CutCopyMode = $false

$spotValues = @{
    2  = 100.91
    3  = 95.73
    4  = 89.09999999999999
    5  = 66.89
    6  = 48.68
    7  = 60.43
    8  = 86.3
    9  = 90.47
    10 = 91.42
    11 = 79.29000000000001
    12 = 33.45
    13 = 19.99
    14 = 40.5
    15 = 17.59
    16 = 11.9
    17 = 16.49
    18 = 39.27
    19 = 49.39
    20 = 70.20999999999999
    21 = 99.98999999999999
    22 = 110.04
    23 = 100.16
    24 = 111
    25 = 93.65000000000001
}

foreach ($row in $spotValues.Keys) {
    $wsSpot.Range("AA$row").Value = $spotValues[$row]
}

# --- Sheet "Gaz": add row 24 (2025-07-08) ---
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A24").NumberFormat = "@"
$wsGaz.Range("A24").Value = "2025-07-08"
$wsGaz.Range("A24").Style = $wsGaz.Range("A23").Style
$wsGaz.Range("B24").Value = 33.85

# --- Sheet "CO2": add row 24 (2025-07-08) ---
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A24").NumberFormat = "@"
$wsCo2.Range("A24").Value = "2025-07-08"
$wsCo2.Range("A24").Style = $wsCo2.Range("A23").Style
$wsCo2.Range("B24").Value = 69.7
